$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "304.00"
Set-TextValue "E2" "4.23%"
Set-TextValue "G2" "17"
Set-TextValue "D3" "35.94"
Set-TextValue "G3" "17"
Set-TextValue "D4" "5.089"
Set-TextValue "E4" "2.80%"
Set-TextValue "G4" "17"
Set-TextValue "D5" "0.07871"
Set-TextValue "E5" "5.62%"
Set-TextValue "G5" "17"
Set-TextValue "D6" "2.285"
Set-TextValue "E6" "3.69%"
Set-TextValue "G6" "17"
Set-TextValue "D7" "8.084"
Set-TextValue "E7" "4.49%"
Set-TextValue "G7" "17"
Set-TextValue "D8" "4.008"
Set-TextValue "E8" "6.91%"
Set-TextValue "G8" "17"
Set-TextValue "D9" "0.9250"
Set-TextValue "E9" "0.50%"
Set-TextValue "G9" "17"
Set-TextValue "D10" "0.1003"
Set-TextValue "E10" "6.53%"
Set-TextValue "G10" "17"
Set-TextValue "D11" "0.1828"
Set-TextValue "E11" "6.22%"
Set-TextValue "G11" "17"
Set-TextValue "D12" "0.08632"
Set-TextValue "E12" "3.67%"
Set-TextValue "G12" "17"
Set-TextValue "D13" "0.03403"
Set-TextValue "E13" "6.92%"
Set-TextValue "G13" "17"
Set-TextValue "D14" "0.09900"
Set-TextValue "E14" "-0.28%"
Set-TextValue "G14" "17"
Set-TextValue "D15" "0.001477"
Set-TextValue "E15" "-1.20%"
Set-TextValue "G15" "17"
Set-TextValue "D16" "0.04662"
Set-TextValue "E16" "3.74%"
Set-TextValue "G16" "17"
Set-TextValue "D17" "0.005595"
Set-TextValue "E17" "-1.64%"
Set-TextValue "G17" "17"
Set-TextValue "D18" "3.484"
Set-TextValue "E18" "0.24%"
Set-TextValue "G18" "17"
Set-TextValue "E19" "-1.57%"
Set-TextValue "G19" "17"
Set-TextValue "D20" "0.3435"
Set-TextValue "E20" "3.21%"
Set-TextValue "G20" "17"
Set-TextValue "E21" "1.44%"
Set-TextValue "G21" "17"
Set-TextValue "D22" "4.552"
Set-TextValue "E22" "9.56%"
Set-TextValue "G22" "17"
Set-TextValue "E23" "5.57%"
Set-TextValue "G23" "17"
Set-TextValue "D24" "0.001241"
Set-TextValue "E24" "2.17%"
Set-TextValue "G24" "17"
Set-TextValue "D25" "0.004491"
Set-TextValue "E25" "5.43%"
Set-TextValue "G25" "17"
Set-TextValue "D26" "0.0001299"
Set-TextValue "E26" "0.20%"
Set-TextValue "G26" "17"
Set-TextValue "D27" "0.0002796"
Set-TextValue "E27" "-17.28%"
Set-TextValue "G27" "17"
Set-TextValue "G28" "17"
Set-TextValue "G29" "17"
Set-TextValue "G30" "17"
Set-TextValue "G31" "17"
Set-TextValue "G32" "17"
Set-TextValue "G33" "17"
Set-TextValue "G34" "17"
Set-TextValue "G35" "17"
Set-TextValue "G36" "17"
Set-TextValue "G37" "17"
Set-TextValue "G38" "17"
Set-TextValue "D39" "0.01758"
Set-TextValue "E39" "8.71%"
Set-TextValue "G39" "17"
Set-TextValue "D40" "0.04688"
Set-TextValue "E40" "2.69%"
Set-TextValue "G40" "17"
Set-TextValue "D41" "0.007886"
Set-TextValue "E41" "6.52%"
Set-TextValue "G41" "17"
Set-TextValue "D42" "0.1413"
Set-TextValue "E42" "4.37%"
Set-TextValue "G42" "17"
Set-TextValue "D43" "0.008786"
Set-TextValue "G43" "17"
Set-TextValue "D44" "0.002217"
Set-TextValue "E44" "2.91%"
Set-TextValue "G44" "17"
Set-TextValue "D45" "0.009200"
Set-TextValue "E45" "-4.69%"
Set-TextValue "G45" "17"
Set-TextValue "D46" "0.00006009"
Set-TextValue "E46" "-1.24%"
Set-TextValue "G46" "17"
Set-TextValue "E47" "0.15%"
Set-TextValue "G47" "17"
Set-TextValue "D48" "3.898"
Set-TextValue "E48" "48.62%"
Set-TextValue "G48" "17"
Set-TextValue "E49" "34.92%"
Set-TextValue "G49" "17"
Set-TextValue "E50" "0.15%"
Set-TextValue "G50" "17"
Set-TextValue "D51" "0.0001998"
Set-TextValue "E51" "0.15%"
Set-TextValue "G51" "17"

Write-Output "Updated symbol list values"
